{"js": "// Revert the multi-industry template text back to the Finance / Core\n// Banking Systems Modernization wording.\n//\n// Each replacement below targets one exact, whole-string occurrence that\n// lives inside a single <w:t> run, so a plain text search-and-replace\n// (via Range.search + Range.insertText(\"Replace\")) is sufficient and keeps\n// all existing run/paragraph formatting untouched.\n\nconst replacements = [\n  [\n    \"ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING\",\n    \"FINANCE - CORE BANKING SYSTEM MODERNIZATION\",\n  ],\n  [\n    \"Industry: Finance and Machine Learning\",\n    \"Industry: Finance and Banking Operations\",\n  ],\n  [\n    \"This project proposal outlines a strategic Finance Implementation initiative for Finance and Machine Learning to achieve Digital transformation through intelligent automation and predictive analytics. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization.\",\n    \"This project proposal outlines a strategic Finance Implementation initiative for Finance and Banking Operations to achieve Digital transformation through intelligent automation and predictive analytics. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization.\",\n  ],\n  [\n    \"Project Name: Finance and Machine Learning Implementation Initiative\",\n    \"Project Name: Finance and Banking Operations Implementation Initiative\",\n  ],\n  [\n    \"Industry Focus: Finance and Machine Learning\",\n    \"Industry Focus: Finance and Banking Operations\",\n  ],\n  [\n    \"\u2022 Regulatory compliance\",\n    \"\u2022 Finance compliance\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const range of found.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Revert the multi-industry template text back to the Finance / Core\n# Banking Systems Modernization wording.\n#\n# Each replacement below targets one exact, whole-string occurrence, so a\n# plain Find/Replace (wdReplaceAll) against the whole document Range is\n# sufficient and leaves every other run/paragraph untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace did not match any text: $findText\"\n    }\n}\n\nReplace-Text \"ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING\" \"FINANCE - CORE BANKING SYSTEM MODERNIZATION\"\n\nReplace-Text \"Industry: Finance and Machine Learning\" \"Industry: Finance and Banking Operations\"\n\nReplace-Text \"This project proposal outlines a strategic Finance Implementation initiative for Finance and Machine Learning to achieve Digital transformation through intelligent automation and predictive analytics. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization.\" \"This project proposal outlines a strategic Finance Implementation initiative for Finance and Banking Operations to achieve Digital transformation through intelligent automation and predictive analytics. The proposed solution addresses critical business challenges while delivering measurable value through Operational efficiency, Predictive maintenance, Customer personalization.\"\n\nReplace-Text \"Project Name: Finance and Machine Learning Implementation Initiative\" \"Project Name: Finance and Banking Operations Implementation Initiative\"\n\nReplace-Text \"Industry Focus: Finance and Machine Learning\" \"Industry Focus: Finance and Banking Operations\"\n\nReplace-Text \"Regulatory compliance\" \"Finance compliance\"\n"}
